$d = $word.ActiveDocument

# Remove the whole paragraph "{% load docx_tags %}" (it spans 3 runs: "{% load docx_",
# "tags", " %}"). Walk backwards so deleting doesn't perturb not-yet-visited indices.
$paras = $d.Paragraphs
for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "{% load docx_tags %}`r") {
        $p.Range.Delete()
    }
}

# Flip the Normal style's "overflowPunct" from false to true
# (exposed on the Word object model as ParagraphFormat.HangingPunctuation).
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.HangingPunctuation = $true
